$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows before row 5 (old row4 shifts down); we build on top of existing row4 ---
$ws.Rows("5:7").Insert() | Out-Null

# --- Update row 3 values (now reflects test3.plusAncs.fa data) ---
$ws.Range("B3").Value = "test_data/test_rareVariants/test3/test3.plusAncs.fa"
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 2
$ws.Range("J3").Value = "pop1_anc"
$ws.Range("K3").Value = "pop2_anc"

# --- Update row 4 values (now reflects test3a.fa data) ---
$ws.Range("B4").Value = "test_data/test_rareVariants/test3/test3a.fa"
$ws.Range("M4").Value = "N.A."
$ws.Range("O4").Value = 0
$ws.Range("T4").Value = ""

# --- Populate new row 5: test3a.plusAncs.fa ---
$ws.Range("A5").Value = "pop1_vs_pop2_Dn3"
$ws.Range("B5").Value = "test_data/test_rareVariants/test3/test3a.plusAncs.fa"
$ws.Range("C5").Value = 13
$ws.Range("D5").Value = 11
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = 20
$ws.Range("J5").Value = "pop1_anc"
$ws.Range("K5").Value = "pop2_anc"
$ws.Range("L5").Value = "N.A."
$ws.Range("M5").Value = "N.A."
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = "N.A."
$ws.Range("S5").Value = "N.A."

# --- Populate new row 6: test3b.fa ---
$ws.Range("A6").Value = "pop1_vs_pop2_Dn4"
$ws.Range("B6").Value = "test_data/test_rareVariants/test3/test3b.fa"
$ws.Range("C6").Value = 11
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = 60
$ws.Range("I6").Value = 20
$ws.Range("J6").Value = "pop1_seq01"
$ws.Range("K6").Value = "pop2_seq01"
$ws.Range("L6").Value = "N.A."
$ws.Range("M6").Value = "1"
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = "N.A."
$ws.Range("S6").Value = "N.A."
$ws.Range("T6").Value = "not signif"

# --- Populate new row 7: test3b.plusAncs.fa ---
$ws.Range("A7").Value = "pop1_vs_pop2_Dn5"
$ws.Range("B7").Value = "test_data/test_rareVariants/test3/test3b.plusAncs.fa"
$ws.Range("C7").Value = 13
$ws.Range("D7").Value = 11
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = 60
$ws.Range("I7").Value = 20
$ws.Range("J7").Value = "pop1_anc"
$ws.Range("K7").Value = "pop2_anc"
$ws.Range("L7").Value = "N.A."
$ws.Range("M7").Value = "1"
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = "N.A."
$ws.Range("S7").Value = "N.A."
$ws.Range("T7").Value = "not signif"

# --- Apply consistent cell formatting from row 4 to the new rows ---
$ws.Range("A4:T4").Copy() | Out-Null
$ws.Range("A5:T7").PasteSpecial(-4122) | Out-Null

# Re-apply the values after formatting paste (PasteSpecial(xlPasteFormats) should not touch
# values, but keep assignment order safe by re-setting key text/number cells once more)
$ws.Range("A5").Value = "pop1_vs_pop2_Dn3"
$ws.Range("B5").Value = "test_data/test_rareVariants/test3/test3a.plusAncs.fa"
$ws.Range("A6").Value = "pop1_vs_pop2_Dn4"
$ws.Range("B6").Value = "test_data/test_rareVariants/test3/test3b.fa"
$ws.Range("A7").Value = "pop1_vs_pop2_Dn5"
$ws.Range("B7").Value = "test_data/test_rareVariants/test3/test3b.plusAncs.fa"

# --- Column B width change ---
$ws.Columns("B").ColumnWidth = 62.005
